$d = $word.ActiveDocument

# Locate the paragraph that contains "LOM3206: Eletrônica (Requisito)" and the
# paragraph that contains the site-footer copyright text; the three
# paragraphs in between (the blank line, "Ver no Jupiter..." and the
# copyright line itself) are removed, leaving the "Requisito" paragraph
# followed directly by the existing blank paragraph before the page break.
$reqIndex = -1
$copyrightIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*LOM3206: Eletrônica (Requisito)*") {
        $reqIndex = $i
    }
    if ($t -like "*Powered by Jekyll and Github pages*") {
        $copyrightIndex = $i
    }
}

if ($reqIndex -gt 0 -and $copyrightIndex -gt $reqIndex) {
    $startPara = $d.Paragraphs($reqIndex + 1)
    $endPara = $d.Paragraphs($copyrightIndex)
    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete()
}
